$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1809045226130653
$ws.Range("C2").Value = 0.6080402010050251
$ws.Range("J2").Value = 0.005025125628140704
$ws.Range("P2").Value = 0.1407035175879397
$ws.Range("S2").Value = 0.06532663316582915
$ws.Range("B3").Value = 0.01646090534979424
$ws.Range("C3").Value = 0.0205761316872428
$ws.Range("J3").Value = 0.02469135802469136
$ws.Range("P3").Value = 0.7407407407407407
$ws.Range("S3").Value = 0.1975308641975309
$ws.Range("J4").Value = 0.05882352941176471
$ws.Range("P4").Value = 0.6078431372549019
$ws.Range("S4").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.04716981132075472
$ws.Range("D6").Value = 0.004716981132075472
$ws.Range("E6").Value = 0.004716981132075472
$ws.Range("F6").Value = 0.04716981132075472
$ws.Range("J6").Value = 0.2547169811320755
$ws.Range("O6").Value = 0.02830188679245283
$ws.Range("Q6").Value = 0.1415094339622641
$ws.Range("R6").Value = 0.1037735849056604
$ws.Range("S6").Value = 0.3679245283018868
$ws.Range("B7").Value = 0.1369294605809129
$ws.Range("D7").Value = 0.03734439834024896
$ws.Range("F7").Value = 0.05394190871369295
$ws.Range("J7").Value = 0.1618257261410788
$ws.Range("O7").Value = 0.004149377593360996
$ws.Range("Q7").Value = 0.1037344398340249
$ws.Range("R7").Value = 0.1037344398340249
$ws.Range("S7").Value = 0.3983402489626556
$ws.Range("B8").Value = 0.1206896551724138
$ws.Range("D8").Value = 0.02155172413793104
$ws.Range("E8").Value = 0.002155172413793103
$ws.Range("F8").Value = 0.04956896551724138
$ws.Range("J8").Value = 0.146551724137931
$ws.Range("O8").Value = 0.02586206896551724
$ws.Range("Q8").Value = 0.1551724137931035
$ws.Range("R8").Value = 0.08836206896551724
$ws.Range("S8").Value = 0.3900862068965517
$ws.Range("B9").Value = 0.1134020618556701
$ws.Range("D9").Value = 0.02061855670103093
$ws.Range("F9").Value = 0.06701030927835051
$ws.Range("J9").Value = 0.1237113402061856
$ws.Range("O9").Value = 0.02577319587628866
$ws.Range("Q9").Value = 0.1597938144329897
$ws.Range("R9").Value = 0.09793814432989691
$ws.Range("S9").Value = 0.3917525773195876
$ws.Range("B10").Value = 0.1256345177664974
$ws.Range("D10").Value = 0.01776649746192894
$ws.Range("E10").Value = 0.001269035532994924
$ws.Range("F10").Value = 0.0532994923857868
$ws.Range("J10").Value = 0.1440355329949239
$ws.Range("O10").Value = 0.01522842639593909
$ws.Range("Q10").Value = 0.2068527918781726
$ws.Range("R10").Value = 0.09137055837563451
$ws.Range("S10").Value = 0.3445431472081218
$ws.Range("F11").Value = 0.002583979328165375
$ws.Range("G11").Value = 0.1498708010335917
$ws.Range("J11").Value = 0.09819121447028424
$ws.Range("K11").Value = 0.2093023255813954
$ws.Range("L11").Value = 0.5193798449612403
$ws.Range("S11").Value = 0.020671834625323
$ws.Range("G12").Value = 0.7081339712918661
$ws.Range("J12").Value = 0.2105263157894737
$ws.Range("K12").Value = 0.009569377990430622
$ws.Range("L12").Value = 0.03827751196172249
$ws.Range("S12").Value = 0.03349282296650718
$ws.Range("G13").Value = 0.7454545454545455
$ws.Range("J13").Value = 0.2363636363636364
$ws.Range("S13").Value = 0.01818181818181818
$ws.Range("F15").Value = 0.02469135802469136
$ws.Range("H15").Value = 0.168724279835391
$ws.Range("I15").Value = 0.06584362139917696
$ws.Range("J15").Value = 0.345679012345679
$ws.Range("K15").Value = 0.05761316872427984
$ws.Range("M15").Value = 0.00823045267489712
$ws.Range("O15").Value = 0.05349794238683128
$ws.Range("S15").Value = 0.2757201646090535
$ws.Range("F16").Value = 0.01509433962264151
$ws.Range("H16").Value = 0.1471698113207547
$ws.Range("I16").Value = 0.07169811320754717
$ws.Range("J16").Value = 0.4566037735849057
$ws.Range("K16").Value = 0.1094339622641509
$ws.Range("M16").Value = 0.03018867924528302
$ws.Range("N16").Value = 0.003773584905660377
$ws.Range("O16").Value = 0.04150943396226415
$ws.Range("S16").Value = 0.1245283018867925
$ws.Range("F17").Value = 0.01244813278008299
$ws.Range("H17").Value = 0.1452282157676349
$ws.Range("I17").Value = 0.06016597510373444
$ws.Range("J17").Value = 0.454356846473029
$ws.Range("K17").Value = 0.0975103734439834
$ws.Range("M17").Value = 0.01659751037344398
$ws.Range("O17").Value = 0.07053941908713693
$ws.Range("S17").Value = 0.1431535269709543
$ws.Range("F18").Value = 0.008032128514056224
$ws.Range("H18").Value = 0.1405622489959839
$ws.Range("I18").Value = 0.1164658634538153
$ws.Range("J18").Value = 0.4216867469879518
$ws.Range("K18").Value = 0.1285140562248996
$ws.Range("M18").Value = 0.01606425702811245
$ws.Range("N18").Value = 0.004016064257028112
$ws.Range("O18").Value = 0.05622489959839357
$ws.Range("S18").Value = 0.108433734939759
$ws.Range("F19").Value = 0.01339915373765867
$ws.Range("H19").Value = 0.2002820874471086
$ws.Range("I19").Value = 0.07052186177715092
$ws.Range("J19").Value = 0.3878702397743301
$ws.Range("K19").Value = 0.1234132581100141
$ws.Range("M19").Value = 0.02397743300423131
$ws.Range("O19").Value = 0.07052186177715092
$ws.Range("S19").Value = 0.1100141043723554
